$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "Comments" column (G), pushing it to H.
$ws.Range("G1").EntireColumn.Insert()

# Final column widths.
$ws.Columns("E").ColumnWidth = 15.6328125
$ws.Columns("F").ColumnWidth = 13.453125
$ws.Columns("G").ColumnWidth = 13.453125
$ws.Columns("H").ColumnWidth = 31

# Header row (row 3): rename existing headers, add the new one.
$ws.Range("E3").Value = "Structure"
$ws.Range("F3").Value = "Model"
$ws.Range("G3").Value = "Error message"

# "Structure"/Model column (E) - green-filled cells with the JSON model name.
$ws.Range("E4").Value = "ListUsers"
$ws.Range("E5").Value = "SingleUser"
$ws.Range("E7").Value = "ListResources"
$ws.Range("E8").Value = "SingleResource"
$ws.Range("E10").Value = "JobUserCreated"
$ws.Range("E11").Value = "JobUserUpdated"
$ws.Range("E12").Value = "JobUserUpdated"
$ws.Range("E14").Value = "RegUser"
$ws.Range("E16").Value = "RegUser"
$ws.Range("E18").Value = "ListUsers"

$greenCells = @("E4","E5","E7","E8","E10","E11","E12","E14","E16","E18")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Interior.Color = 5296274
}

# Column F "Check Exact" markers.
$ws.Range("F5").Value = "Check Exact"
$ws.Range("F8").Value = "Check Exact"

# Column G "Check Error" markers.
$ws.Range("G6").Value = "Check Error"
$ws.Range("G9").Value = "Check Error"
$ws.Range("G15").Value = "Check Error"
$ws.Range("G17").Value = "Check Error"

# Update the view: scroll back to show row 1, select G18.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("G18").Select()
